$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Paragraph")
$ws.Activate()

# Column C ("Revised Format") for rows 71-83 used to just repeat the
# "Sample_code" value copied from column B (Original Format) - i.e. no
# real rename had been recorded yet. Give it its own revised value,
# "SampleCode", instead of aliasing the original format's text.
for ($r = 71; $r -le 83; $r++) {
    $ws.Cells.Item($r, 3).Value = "SampleCode"
}

# Mirror the author's scrolled/selected view state: window scrolled a
# little further down and the active cell moved to D78.
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D78").Select()
